$d = $word.ActiveDocument

function Assert-ParaText($index, $expected) {
    # Paragraph.Range.Text includes the trailing paragraph mark.
    $actual = $d.Paragraphs.Item($index).Range.Text
    $expectedWithMark = $expected + "`r"
    if ($actual -ne $expectedWithMark) {
        throw "Paragraph $index expected [$expectedWithMark] but found [$actual]"
    }
}

# --- 1) "Javascript" -> "REACT" + "js" (two separate runs, same formatting) ---
Assert-ParaText 8 "Javascript"
$p = $d.Paragraphs.Item(8)
$r = $p.Range
$rText = $d.Range($r.Start, $r.End - 1)
$rText.Text = "REACT"
$p2 = $d.Paragraphs.Item(8)
$rEnd = $d.Range($p2.Range.End - 1, $p2.Range.End - 1)
$rEnd.InsertAfter("js")
# Toggling formatting on the inserted text and back forces Word to keep it
# as its own run (same visible rPr) instead of silently merging it with
# the preceding "REACT" run.
$rEnd.Font.Bold = 1
$rEnd.Font.Bold = 0
Assert-ParaText 8 "REACTjs"

# --- 2) Remove the old standalone "REACT" paragraph (now paragraph 9) ---
Assert-ParaText 9 "REACT"
$d.Paragraphs.Item(9).Range.Delete()

# --- 3) "JSX (XML/HTML)" + tab + "?" -> just "JSX (XML/HTML)" ---
Assert-ParaText 9 "JSX (XML/HTML)`t?"
$null = $d.Content.Find.Execute("JSX (XML/HTML)`t?", $true, $false, $false, $false, $false, $true, 1, $false, "JSX (XML/HTML)", 2)
Assert-ParaText 9 "JSX (XML/HTML)"

# --- 4) Remove "Boosts JS" paragraph and the empty paragraph after it ---
Assert-ParaText 10 "Boosts JS"
$d.Paragraphs.Item(10).Range.Delete()
Assert-ParaText 10 ""
$d.Paragraphs.Item(10).Range.Delete()
Assert-ParaText 10 "Backend"

# --- 5) Insert REDUX / Web API / Node / MongoDB paragraphs before "Backend" ---
$backend = $d.Paragraphs.Item(10)
$insertPoint = $d.Range($backend.Range.Start, $backend.Range.Start)
$insertPoint.InsertBefore("REDUX`rWeb API`rNode`rMongoDB`r")

# New paragraphs inherit "Backend"'s list level on insertion; fix them up to
# match the intended outline: REDUX/Node/MongoDB one level deeper than
# "Web API" (which sits at the same level as "Frontend"/"Backend").
$d.Paragraphs.Item(10).Range.ListFormat.ListLevelNumber = 3
$d.Paragraphs.Item(11).Range.ListFormat.ListLevelNumber = 2
$d.Paragraphs.Item(12).Range.ListFormat.ListLevelNumber = 3
$d.Paragraphs.Item(13).Range.ListFormat.ListLevelNumber = 3

Assert-ParaText 10 "REDUX"
Assert-ParaText 11 "Web API"
Assert-ParaText 12 "Node"
Assert-ParaText 13 "MongoDB"
Assert-ParaText 14 "Backend"
